$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE (J2): "001" -> "002"
# Force text storage (avoid "002" being coerced to the number 2), then
# restore the cell's style so no explicit number-format style sticks to it.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").Style = "Normal"

# REPORT_TYPE_CODE (K2) stays "001" - unchanged.

# NOTICE_DATE / REPORT_DATE
$ws.Range("M2").Value = "2020-12-18 00:00:00"
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Financial figures
$ws.Range("O2").Value = 3138321859.39
$ws.Range("P2").Value = 1249765848.17
$ws.Range("Q2").Value = 244120322.42
$ws.Range("S2").Value = 465713337.53
$ws.Range("U2").Value = 664743052.6799999
$ws.Range("W2").Value = 1881273163.17
$ws.Range("X2").Value = 526469615.41

# ADVANCE_RECEIVABLES (Z2) becomes blank (was numeric)
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = ""
$ws.Range("Z2").Style = "Normal"

$ws.Range("AB2").Value = 1257048696.22
$ws.Range("AF2").Value = 88.5481090016
$ws.Range("AG2").Value = 59.9451951539
